$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Custom factor" feature: the TYPE column value for the JPY=X row changes
# from "Macro" to "Custom".
$ws.Range("B2").Value = "Custom"

# Simple UI tweak: the active selection moved to D3.
$ws.Range("D3").Select() | Out-Null
